$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as text (matching original formatting)
$dCells = "D2,D3,D5,D7,D9,D10,D12,D14,D15,D16,D17,D18,D19,D21,D22,D24,D26,D29,D30,D31,D32,D34,D36,D39,D41,D42,D44,D46,D47,D49,D50,D51".Split(",")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '42.878.20'
$ws.Range("E2").Value = '  -0.24%  '

$ws.Range("D3").Value = '2.280.40'
$ws.Range("E3").Value = '  -0.31%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '250.59'
$ws.Range("E5").Value = '  -0.70%  '

$ws.Range("E6").Value = '  -1.52%  '

$ws.Range("D7").Value = '79.32'
$ws.Range("E7").Value = '  +9.37%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("D9").Value = '0.638'
$ws.Range("E9").Value = '  -3.77%  '

$ws.Range("D10").Value = '41.50'
$ws.Range("E10").Value = '  +5.74%  '

$ws.Range("E11").Value = '  -1.60%  '

$ws.Range("D12").Value = '7.37'
$ws.Range("E12").Value = '  +0.46%  '

$ws.Range("E13").Value = '  -1.29%  '

$ws.Range("D14").Value = '2.620.12'
$ws.Range("E14").Value = '  -0.29%  '

$ws.Range("D15").Value = '15.22'
$ws.Range("E15").Value = '  +0.57%  '

$ws.Range("D16").Value = '0.869'
$ws.Range("E16").Value = '  -2.87%  '

$ws.Range("D17").Value = '2.276.13'
$ws.Range("E17").Value = '  -0.10%  '

$ws.Range("D18").Value = '42.784.94'
$ws.Range("E18").Value = '  -0.28%  '

$ws.Range("D19").Value = '0.0₃0996'
$ws.Range("E19").Value = '  -1.96%  '

$ws.Range("E20").Value = '  -2.44%  '

$ws.Range("D21").Value = '72.25'
$ws.Range("E21").Value = '  -2.03%  '

$ws.Range("D22").Value = '234.26'
$ws.Range("E22").Value = '  -1.55%  '

$ws.Range("E23").Value = '  +1.03%  '

$ws.Range("D24").Value = '3.78'
$ws.Range("E24").Value = '  -2.44%  '

$ws.Range("E25").Value = '  -0.06%  '

$ws.Range("D26").Value = '11.39'
$ws.Range("E26").Value = '  -2.79%  '

$ws.Range("E27").Value = '  -4.70%  '

$ws.Range("E28").Value = '  +2.11%  '

$ws.Range("D29").Value = '169.08'
$ws.Range("E29").Value = '  +0.33%  '

$ws.Range("D30").Value = '20.91'
$ws.Range("E30").Value = '  -1.44%  '

$ws.Range("D31").Value = '6.66'
$ws.Range("E31").Value = '  +5.47%  '

$ws.Range("D32").Value = '0.0853'
$ws.Range("E32").Value = '  +4.60%  '

$ws.Range("E33").Value = '  -5.19%  '

$ws.Range("D34").Value = '30.46'
$ws.Range("E34").Value = '  -3.42%  '

$ws.Range("E35").Value = '  +1.38%  '

$ws.Range("D36").Value = '4.57'
$ws.Range("E36").Value = '  -4.65%  '

$ws.Range("E37").Value = '  -0.75%  '

$ws.Range("E38").Value = '  -2.60%  '

$ws.Range("D39").Value = '13.57'
$ws.Range("E39").Value = '  +2.39%  '

$ws.Range("E40").Value = '  -2.66%  '

$ws.Range("D41").Value = '5.95'
$ws.Range("E41").Value = '  -1.91%  '

$ws.Range("D42").Value = '115.53'
$ws.Range("E42").Value = '  +18.49%  '

$ws.Range("E43").Value = '  -1.91%  '

$ws.Range("D44").Value = '61.48'
$ws.Range("E44").Value = '  -0.81%  '

$ws.Range("E45").Value = '  -3.33%  '

$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").Value = '0.102'
$ws.Range("E46").Value = '  -2.68%  '

$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").Value = '4.69'
$ws.Range("E47").Value = '  -6.09%  '

$ws.Range("E48").Value = '  -0.05%  '

$ws.Range("D49").Value = '1.15'
$ws.Range("E49").Value = '  -4.18%  '

$ws.Range("D50").Value = '1.18'
$ws.Range("E50").Value = '  -1.92%  '

$ws.Range("D51").Value = '4.31'
$ws.Range("E51").Value = '  -1.07%  '
